$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell text content per the shared-string reshuffle/edit described in the diff.
$ws.Range("A2").Value = 'Matériel'
$ws.Range("C3").Value = 'Icône pas claire ou pas assez visible, textes trop longs ou pas assez vulgarisés, îcones ou raccourcis qui sortent des standards'
$ws.Range("B5").Value = 'Perturbation cognitive au cours d’une tâche. Vous êtes déconcentré par un événement de l’interface,vous avez oublié ce que vous étiez en train de faire.'
$ws.Range("B6").Value = 'Cliquer ou "taper" à côté de l''objet visé : un bouton, un lien, un email, une cible, une touche de clavier etc.Le clic peut soit se produire "dans le vide" (sans conséquence), soit activer un autre objet ou fenêtre non visé(e).'
$ws.Range("C7").Value = 'Oublier la pièce jointe d''un mail, oublier d''enregistrer avant de compiler ou push sur un git,oublier de changer de calque avant de modifier un dessin'
$ws.Range("C8").Value = 'Utilisation du mauvais raccourci, habitude de disposition des icônes qui mène à des erreurs de clics lors d''un changement,fermer trop vite ses notifications'
$ws.Range("B10").Value = 'L''interface change juste avant un clic ou l''appui sur une touche : une pop-up apparait, une liste se met à jour, une autre applicationprend le focus, etc. Cela a pour conséquence que ce clic se produit "dans le vide" ou sur un objet non désiré. '
$ws.Range("C10").Value = 'Mise à jour de la liste d''autocomplétion pendant l''entrée d''un texte, on sélectionne l''élément qui vient tout juste de changer sanspouvoir réagir. Clic "capturé" au dernier moment par une popup ou une notification. Boîte de dialogue ou raccourcis disparaît juste avant le clic.'

# Update the active view: select C3 (and scroll so column C is the leftmost visible column).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C3").Select()
